$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting the existing rows 29-32 down to 30-33
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly price entry
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C29").Value = "Ñuble"
$ws.Range("D29").Value = 44841
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = 100112037
$ws.Range("G29").Value = "Cebollín"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 7500
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 7750
$ws.Range("N29").Value = "$/docena de atados"
$ws.Range("O29").Value = "Provincia de Diguillín"
$ws.Range("P29").Value = 2583
$ws.Range("Q29").Value = 3
$ws.Range("R29").Value = "Hortaliza"
